$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G3").ClearContents()
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "23/7/2020  9:00:00 am"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "`$123.00"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "98.76%"
$excel.CalculateFullRebuild()
